# Update the NATMI LR-pair table with the new TPM-based numbers and add the
# "Inflammatory-Mac" sending cluster rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: FAPs -> Ntng1/Lrrc4c -> FAPs (values refreshed) -------------
$ws.Range("G2").Value = 0.06050633333333334
$ws.Range("I2").Value = 0.6173275744796626
$ws.Range("J2").Value = 0.6173275744796626
$ws.Range("O2").Value = 0.9540047954114494
$ws.Range("P2").Value = 0.9540047954114494
$ws.Range("S2").Value = 0.5889334663933168
$ws.Range("T2").Value = 0.5889334663933168

# --- Row 3: FAPs -> Ntng1/Lrrc4c -> MuSCs (values refreshed) ------------
$ws.Range("G3").Value = 0.06050633333333334
$ws.Range("I3").Value = 0.6173275744796626
$ws.Range("J3").Value = 0.6173275744796626
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01076166666666667
$ws.Range("N3").Value = 0.032285
$ws.Range("O3").Value = 0.04599520458855057
$ws.Range("P3").Value = 0.04599520458855057
$ws.Range("Q3").Value = 0.0006511489905555556
$ws.Range("R3").Value = 0.005860340915000001
$ws.Range("S3").Value = 0.02839410808634577
$ws.Range("T3").Value = 0.02839410808634577

# --- Row 4: becomes Inflammatory-Mac -> Ntng1/Lrrc4c -> FAPs ------------
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.0006813333333333333
$ws.Range("H4").Value = 0.002044
$ws.Range("I4").Value = 0.006951435178887225
$ws.Range("J4").Value = 0.006951435178887225
$ws.Range("O4").Value = 0.9540047954114494
$ws.Range("P4").Value = 0.9540047954114494
$ws.Range("Q4").Value = 0.000152081776
$ws.Range("R4").Value = 0.001368735984
$ws.Range("S4").Value = 0.006631702495650259
$ws.Range("T4").Value = 0.006631702495650259

# --- Row 5: becomes Inflammatory-Mac -> Ntng1/Lrrc4c -> MuSCs -----------
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.0006813333333333333
$ws.Range("H5").Value = 0.002044
$ws.Range("I5").Value = 0.006951435178887225
$ws.Range("J5").Value = 0.006951435178887225
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01076166666666667
$ws.Range("N5").Value = 0.032285
$ws.Range("O5").Value = 0.04599520458855057
$ws.Range("P5").Value = 0.04599520458855057
$ws.Range("Q5").Value = 0.000007332282222222222
$ws.Range("R5").Value = 0.00006599053999999999
$ws.Range("S5").Value = 0.0003197326832369655
$ws.Range("T5").Value = 0.0003197326832369655

# --- New row 6: MuSCs -> Ntng1/Lrrc4c -> FAPs ---------------------------
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ntng1"
$ws.Range("C6").Value = "Lrrc4c"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.03682566666666667
$ws.Range("H6").Value = 0.110477
$ws.Range("I6").Value = 0.3757209903414501
$ws.Range("J6").Value = 0.3757209903414501
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.223212
$ws.Range("N6").Value = 0.669636
$ws.Range("O6").Value = 0.9540047954114494
$ws.Range("P6").Value = 0.9540047954114494
$ws.Range("Q6").Value = 0.008219930707999999
$ws.Range("R6").Value = 0.073979376372
$ws.Range("S6").Value = 0.3584396265224823
$ws.Range("T6").Value = 0.3584396265224823

# --- New row 7: MuSCs -> Ntng1/Lrrc4c -> MuSCs --------------------------
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Ntng1"
$ws.Range("C7").Value = "Lrrc4c"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.03682566666666667
$ws.Range("H7").Value = 0.110477
$ws.Range("I7").Value = 0.3757209903414501
$ws.Range("J7").Value = 0.3757209903414501
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01076166666666667
$ws.Range("N7").Value = 0.032285
$ws.Range("O7").Value = 0.04599520458855057
$ws.Range("P7").Value = 0.04599520458855057
$ws.Range("Q7").Value = 0.0003963055494444445
$ws.Range("R7").Value = 0.003566749945
$ws.Range("S7").Value = 0.01728136381896783
$ws.Range("T7").Value = 0.01728136381896783
